$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 entirely (even_MAG-GUT88709.fa row)
$ws.Rows.Item(3).Delete()

# Delete column C ("max") so that D (prediction) shifts to C, E (rejection-f) shifts to D
$ws.Columns.Item(3).Delete()

# Update row 2 values to match new data
$ws.Range("B2").Value = 10990.46771063232
$ws.Range("C2").Value = "o__Fusobacteriales"
$ws.Range("D2").Value = "o__Fusobacteriales"
